$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.543301582336426
$ws.Range("B1").Value = 3.209999561309814
$ws.Range("C1").Value = 2.938354253768921
$ws.Range("D1").Value = 3.335667610168457
$ws.Range("E1").Value = 1.947994112968445
